$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-sort the data table (A4:S61) by Year (A) descending, then Source (S) ascending ---
# This reproduces the author re-sorting the list after adding/editing rows, which moved
# the 2024 entries (rows 5-8) around: rows 5/6 (EBSCOhost) and 7/8 (AIS) swap places.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A4:A61"), 0, 2, $null, 0)
$sortObj.SortFields.Add($ws.Range("S4:S61"), 0, 1, $null, 0)
$sortObj.SetRange($ws.Range("A4:S61"))
$sortObj.Header = 2
$sortObj.MatchCase = $false
$sortObj.Orientation = 1
$sortObj.Apply()

# --- Re-tag a handful of individual rows to different category columns ---
# Row 10: Team Composition (N) -> Cognition (L)
$ws.Range("N10").Value = $null
$ws.Range("L10").Value = 1

# Row 13: Team Composition (N) -> Trust (J)
$ws.Range("N13").Value = $null
$ws.Range("J13").Value = 1

# Row 30: Team Composition (N) -> Coordination (K)
$ws.Range("N30").Value = $null
$ws.Range("K30").Value = 1

# Row 32: Team Composition (N) -> Coordination (K)
$ws.Range("N32").Value = $null
$ws.Range("K32").Value = 1

# Row 42: Coordination (K) -> Team Composition (N)
$ws.Range("K42").Value = $null
$ws.Range("N42").Value = 1

# Row 54: Team Composition (N) -> Communication (Q)
$ws.Range("N54").Value = $null
$ws.Range("Q54").Value = 1

# --- View state: zoom to 100% and move the selection ---
$excel.ActiveWindow.Zoom = 100
$ws.Range("N63").Select()
